$d = $word.ActiveDocument

# --- Step 1: bump "Five" to "Six" in the intro line ---
$introRng = $d.Content
$introRng.Find.Execute("Five converging forces make this architectural evolution unavoidable.", $true, $false, $false, $false, $false, $true, 1, $false, "Six converging forces make this architectural evolution unavoidable.", 2)
if (-not $introRng.Find.Found) { throw "Could not find the 'Five converging forces' sentence to update." }

# --- Step 2: locate the anchor paragraph = last paragraph of section 5 ("The Quantum Horizon" body) ---
$anchorRng = $d.Content
$anchorRng.Find.Execute("no single point holds the keys to the kingdom.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $anchorRng.Find.Found) { throw "Could not find the end of the Quantum Horizon section." }
$anchor = $anchorRng.Paragraphs(1)

# --- Step 3: capture donor formatting as FormattedText (spanning through the paragraph mark so pPr carries over) ---
# heading donor: an existing numbered section heading ("5. The Quantum Horizon")
$headRng = $d.Content
$headRng.Find.Execute("5. The Quantum Horizon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $headRng.Find.Found) { throw "Could not find the heading donor paragraph." }
$headingDonor = $headRng.Paragraphs(1)
$headingDonorFull = $d.Range($headingDonor.Range.Start, $headingDonor.Next().Range.Start)
$headingFT = $headingDonorFull.FormattedText

# body donor: the anchor paragraph itself (plain body text style used throughout the whitepaper)
$bodyDonorFull = $d.Range($anchor.Range.Start, $anchor.Next().Range.Start)
$bodyFT = $bodyDonorFull.FormattedText

# pull-quote donor: an existing shaded callout paragraph ("Zero-Knowledge Trust is that model.")
$quoteRng = $d.Content
$quoteRng.Find.Execute("Zero-Knowledge Trust is that model.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $quoteRng.Find.Found) { throw "Could not find the pull-quote donor paragraph." }
$quoteDonor = $quoteRng.Paragraphs(1)
$quoteDonorFull = $d.Range($quoteDonor.Range.Start, $quoteDonor.Next().Range.Start)
$quoteFT = $quoteDonorFull.FormattedText

# --- Step 4: insert the new "6. The Human Factor as Superpower" heading right after the anchor ---
$ip = $d.Range($anchor.Range.End, $anchor.Range.End)
$ip.FormattedText = $headingFT
$lastPara = $anchor.Next()
$lastPara.Range.Text = '6. The Human Factor as Superpower'

# --- Step 5: insert the six body paragraphs, each using the body-style donor formatting ---
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'For decades, the security industry has operated under a damaging assumption: people are the weakest link. Training programs, compliance mandates, and ever-more-restrictive access policies all stem from the same premise—that users are liabilities to be managed. This framing is not just demoralizing; it is a confession that the security model requires perfect human behavior to function. When it inevitably fails, it blames the people it was supposed to protect.'
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'Zero-trust architectures, for all their sophistication, perpetuate this dynamic. They verify identity rigorously—then hand over plaintext secrets and hope the verified human handles them correctly. The breach does not happen at the verification step. It happens after, when a trusted administrator misconfigures a policy, a developer pastes a credential into a chat, or a phished employee surrenders an access token. The architecture works perfectly until a person behaves like a person.'
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'Zero-Knowledge Trust reframes the human factor entirely. Rather than treating people as failure points to be constrained, it treats them as the creative, adaptive, resourceful actors they are—and builds an architecture that lets them operate at full capability without risk of exposure. The system does not restrict what users can do; it ensures that their actions cannot produce catastrophic outcomes.'
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'Consider the practical difference. In a zero-trust environment, an employee must navigate a maze of access restrictions, remember which secrets go where, avoid sharing credentials through unapproved channels, and follow rotation schedules they did not design. Every one of these is a friction point that slows work and creates opportunities for error. The model treats the human as an adversary of its own security.'
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'In a Zero-Knowledge Trust environment, the same employee simply works. Secrets are delivered ephemerally to the systems that need them, scoped to the operation, and never materialize in a form that can be copied, shared, or mishandled. There is no credential to paste into a chat because there is no credential the user ever sees. There is no rotation to forget because rotation is cryptographic and automatic. There is no misconfiguration that exposes a vault because the vault cannot be read by the infrastructure that hosts it.'
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $bodyFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'The result is not just better security—it is better work. When people are freed from the cognitive burden of being their own security layer, they move faster, collaborate more openly, and focus their energy on the work that matters. The human factor transforms from the industry’s perennial excuse for failure into its greatest competitive advantage.'

# --- Step 6: insert the closing pull-quote paragraph using the pull-quote donor formatting ---
$ip = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$ip.FormattedText = $quoteFT
$lastPara = $lastPara.Next()
$lastPara.Range.Text = 'The strongest security model is not the one that expects the least of people. It is the one that enables the most—while making exposure architecturally impossible.'

Write-Output "edit complete"